$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cadastro")

# Fix cadastro screen data: correct user name value in A2
$ws.Range("A2").Value = "AngraSouzaaa"

# Update selection to reflect the active cell after the fix
$ws.Range("A2").Select()
